$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$ws3 = $wb.Worksheets.Item("Burndow - Sprint2")
$ws3.Delete()

$excel.DisplayAlerts = $true
